$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update header text: "Model size in bytes" -> "Model size in Mega Bytes"
$ws.Range("D1").Value = "Model size in Mega Bytes"

# 2. Convert "Model size" column values from bytes to Mega Bytes (divide by 1,000,000)
for ($row = 2; $row -le 8; $row++) {
    $cell = $ws.Cells.Item($row, 4)
    $bytes = $cell.Value2
    $cell.Value = $bytes / 1000000
}

# 3. Update the default alignment for column A (Model name) from general to left
$ws.Columns.Item(1).HorizontalAlignment = -4131   # xlLeft
